$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.283.62"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "2.639.50"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.66"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.176"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.64%  "
$ws.Range("D10").Value = "2.639.26"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E14").Value = "  +4.47%  "
$ws.Range("D15").Value = "3.122.41"
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").Value = "72.171.47"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.57"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "2.640.80"
$ws.Range("E18").Value = "  +2.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.95"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.92"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "378.36"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +10.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.07"
$ws.Range("D24").ClearFormats()
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.39"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.01"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("D28").Value = "2.776.84"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "0.0₃0957"
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "526.79"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.14"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.82"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.65"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.33"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("E38").Value = "  -6.05%  "
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.87"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.33"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.32"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.545"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "0.0₆0260"
$ws.Range("E51").Value = "  -4.68%  "
